# Update the EPEX Spot prices workbook with the latest daily data point.
$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add a new day column (BN) with the 18-aug prices ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$spotValues = @(
    76.15000000000001,
    56.54,
    74.12,
    62.61,
    63.28,
    68.15000000000001,
    75.45,
    48.64,
    67.54000000000001,
    78.86,
    61.1,
    45.93,
    20.04,
    14.23,
    18.34,
    54.56,
    69.26000000000001,
    84.53,
    104.07,
    115.06,
    120,
    120.76,
    100.44,
    79.56
)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 66).Value = $spotValues[$i]
}

# Copy the header style (bold, centered, bordered) from BM1 onto the new BN1 header cell.
$wsSpot.Range("BM1").Copy() | Out-Null
$wsSpot.Range("BN1").PasteSpecial(-4122) | Out-Null
$wsSpot.Range("BN1").Value = "18-aug"

# --- Sheet "Gaz": append the two new daily rows ---
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A63").NumberFormat = "@"
$wsGaz.Range("A63").Value = "2025-08-16"
$wsGaz.Range("A63").Style = "Normal"
$wsGaz.Range("B63").Value = 29.925

$wsGaz.Range("A64").NumberFormat = "@"
$wsGaz.Range("A64").Value = "2025-08-17"
$wsGaz.Range("A64").Style = "Normal"
$wsGaz.Range("B64").Value = 29.925

# --- Sheet "CO2": append the two new daily rows ---
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A63").NumberFormat = "@"
$wsCo2.Range("A63").Value = "2025-08-16"
$wsCo2.Range("A63").Style = "Normal"
$wsCo2.Range("B63").Value = 69.95

$wsCo2.Range("A64").NumberFormat = "@"
$wsCo2.Range("A64").Value = "2025-08-17"
$wsCo2.Range("A64").Style = "Normal"
$wsCo2.Range("B64").Value = 69.95
